# Fruta / hortaliza, semanal
# Insert two new price records (dated 2022-08-03 / serial 44776) right after
# the existing row 518 (Naranja - Valencia, 2021-05-14), pushing the rest of
# the data block down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 519 (shifts old rows 519..586 -> 521..588)
$ws.Rows.Item(519).Insert()
$ws.Rows.Item(519).Insert()

# --- New row 519: Naranja / Fukumoto / Primera ---
$ws.Cells.Item(519, 1).Value  = 5
$ws.Cells.Item(519, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(519, 3).Value  = "Maule"
$ws.Cells.Item(519, 4).Value  = 44776
$ws.Cells.Item(519, 5).Value  = 7
$ws.Cells.Item(519, 6).Value  = "Fruta"
$ws.Cells.Item(519, 7).Value  = 100102
$ws.Cells.Item(519, 8).Value  = "Cítricos"
$ws.Cells.Item(519, 9).Value  = 100102005
$ws.Cells.Item(519, 10).Value = "Naranja"
$ws.Cells.Item(519, 11).Value = "Fukumoto"
$ws.Cells.Item(519, 12).Value = "Primera"
$ws.Cells.Item(519, 13).Value = 300
$ws.Cells.Item(519, 14).Value = 6000
$ws.Cells.Item(519, 15).Value = 6000
$ws.Cells.Item(519, 16).Value = 6000
$ws.Cells.Item(519, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(519, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(519, 19).Value = 400
$ws.Cells.Item(519, 20).Value = 15

# --- New row 520: Naranja / Navel Late / Primera ---
$ws.Cells.Item(520, 1).Value  = 5
$ws.Cells.Item(520, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(520, 3).Value  = "Maule"
$ws.Cells.Item(520, 4).Value  = 44776
$ws.Cells.Item(520, 5).Value  = 7
$ws.Cells.Item(520, 6).Value  = "Fruta"
$ws.Cells.Item(520, 7).Value  = 100102
$ws.Cells.Item(520, 8).Value  = "Cítricos"
$ws.Cells.Item(520, 9).Value  = 100102005
$ws.Cells.Item(520, 10).Value = "Naranja"
$ws.Cells.Item(520, 11).Value = "Navel Late"
$ws.Cells.Item(520, 12).Value = "Primera"
$ws.Cells.Item(520, 13).Value = 300
$ws.Cells.Item(520, 14).Value = 5000
$ws.Cells.Item(520, 15).Value = 5000
$ws.Cells.Item(520, 16).Value = 5000
$ws.Cells.Item(520, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(520, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(520, 19).Value = 333
$ws.Cells.Item(520, 20).Value = 15
